# Update the "Clasificación" standings table with the latest matchday
# results. One extra game (PJ) has been played by every team, which
# shifts PG/PP/TD/TP/DT/V1/PTS accordingly for rows 2-9 (players David,
# Pedro, Adonay, Iván, Nico, Nicolás, Richard, Vicente).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - David
$ws.Range("C2").Value = 5.0
$ws.Range("D2").Value = 4.0
$ws.Range("G2").Value = 3.0
$ws.Range("I2").Value = 2.0
$ws.Range("J2").Value = 4.0
$ws.Range("M2").Value = 12.0

# Row 3 - Pedro
$ws.Range("C3").Value = 5.0
$ws.Range("F3").Value = 4.0
$ws.Range("H3").Value = 7.0
$ws.Range("I3").Value = -5.0

# Row 4 - Adonay
$ws.Range("C4").Value = 5.0
$ws.Range("D4").Value = 5.0
$ws.Range("G4").Value = 6.0
$ws.Range("I4").Value = 6.0
$ws.Range("J4").Value = 3.0
$ws.Range("M4").Value = 16.0

# Row 5 - Richard
$ws.Range("C5").Value = 5.0
$ws.Range("F5").Value = 4.0
$ws.Range("H5").Value = 4.0
$ws.Range("I5").Value = -4.0

# Row 6 - Iván
$ws.Range("C6").Value = 5.0
$ws.Range("F6").Value = 2.0
$ws.Range("H6").Value = 3.0
$ws.Range("I6").Value = 4.0

# Row 7 - Nico
$ws.Range("C7").Value = 5.0
$ws.Range("D7").Value = 2.0
$ws.Range("G7").Value = 2.0
$ws.Range("I7").Value = 1.0
$ws.Range("J7").Value = 2.0
$ws.Range("M7").Value = 7.0

# Row 8 - Nicolás
$ws.Range("C8").Value = 5.0
$ws.Range("F8").Value = 4.0
$ws.Range("H8").Value = 6.0
$ws.Range("I8").Value = -5.0

# Row 9 - Vicente
$ws.Range("C9").Value = 5.0
$ws.Range("D9").Value = 3.0
$ws.Range("G9").Value = 3.0
$ws.Range("I9").Value = 1.0
$ws.Range("J9").Value = 3.0
